# Scheduled market-data refresh: update Leve profit calculation sheets (H-N columns)
# across all 8 job worksheets (ALC, ARM, BSM, CRP, CUL, GSM, LTW, WVR).
$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
# Row 8
$ws.Range("H8").Value = 783
$ws.Range("I8").Value = 139.6
$ws.Range("K8").Value = 418.8
$ws.Range("M8").Value = -279.8

# Row 113
$ws.Range("H113").Value = 8386.883
$ws.Range("I113").Value = 12464.333
$ws.Range("J113").Value = 3799.75
$ws.Range("K113").Value = 12464.333
$ws.Range("L113").Value = 3799.75
$ws.Range("M113").Value = -9210.333000000001
$ws.Range("N113").Value = -10307.75

# Row 116
$ws.Range("H116").Value = 6345.8096
$ws.Range("I116").Value = 3743.5715
$ws.Range("J116").Value = 7646.9287
$ws.Range("K116").Value = 3743.5715
$ws.Range("L116").Value = 7646.9287
$ws.Range("M116").Value = -301.5715
$ws.Range("N116").Value = -14530.9287

# Row 132
$ws.Range("H132").Value = 32804242
$ws.Range("I132").Value = 39010124
$ws.Range("J132").Value = 1714.2858
$ws.Range("K132").Value = 117030372
$ws.Range("L132").Value = 5142.857400000001
$ws.Range("M132").Value = -117027842
$ws.Range("N132").Value = -10202.8574

# Row 138
$ws.Range("H138").Value = 1752.5245
$ws.Range("I138").Value = 1126.4
$ws.Range("J138").Value = 2358.4517
$ws.Range("K138").Value = 3379.2
$ws.Range("L138").Value = 7075.355100000001
$ws.Range("M138").Value = 1760.8
$ws.Range("N138").Value = -17355.3551

$ws = $wb.Worksheets.Item("ARM")
# Row 2
$ws.Range("H2").Value = 918.3333
$ws.Range("I2").Value = 870.3333
$ws.Range("J2").Value = 966.3333
$ws.Range("K2").Value = 870.3333
$ws.Range("L2").Value = 966.3333
$ws.Range("M2").Value = -757.3333
$ws.Range("N2").Value = -1192.3333

# Row 45
$ws.Range("H45").Value = 1328.0834
$ws.Range("I45").Value = 1915.6666
$ws.Range("J45").Value = 740.5
$ws.Range("K45").Value = 1915.6666
$ws.Range("L45").Value = 740.5
$ws.Range("M45").Value = -1538.6666
$ws.Range("N45").Value = -1494.5

# Row 61
$ws.Range("H61").Value = 11473.6
$ws.Range("I61").Value = 13279.5
$ws.Range("J61").Value = 4250
$ws.Range("K61").Value = 13279.5
$ws.Range("L61").Value = 4250
$ws.Range("M61").Value = -13067.5
$ws.Range("N61").Value = -4674

# Row 74
$ws.Range("H74").Value = 3408.3845
$ws.Range("I74").Value = 576.5806
$ws.Range("J74").Value = 14381.625
$ws.Range("K74").Value = 576.5806
$ws.Range("L74").Value = 14381.625
$ws.Range("M74").Value = 297.4194
$ws.Range("N74").Value = -16129.625

# Row 77
$ws.Range("H77").Value = 3408.3845
$ws.Range("I77").Value = 576.5806
$ws.Range("J77").Value = 14381.625
$ws.Range("K77").Value = 2882.903
$ws.Range("L77").Value = 71908.125
$ws.Range("M77").Value = 1485.097
$ws.Range("N77").Value = -80644.125

# Row 116
$ws.Range("H116").Value = 918.3333
$ws.Range("I116").Value = 870.3333
$ws.Range("J116").Value = 966.3333
$ws.Range("K116").Value = 870.3333
$ws.Range("L116").Value = 966.3333
$ws.Range("M116").Value = 1423.6667
$ws.Range("N116").Value = -5554.3333

# Row 122
$ws.Range("H122").Value = 75002120
$ws.Range("I122").Value = 133334904
$ws.Range("J122").Value = 2822.8572
$ws.Range("K122").Value = 400004712
$ws.Range("L122").Value = 8468.571599999999
$ws.Range("M122").Value = -400002262
$ws.Range("N122").Value = -13368.5716

# Row 136
$ws.Range("H136").Value = 11473.6
$ws.Range("I136").Value = 13279.5
$ws.Range("J136").Value = 4250
$ws.Range("K136").Value = 39838.5
$ws.Range("L136").Value = 12750
$ws.Range("M136").Value = -37288.5
$ws.Range("N136").Value = -17850

$ws = $wb.Worksheets.Item("BSM")
# Row 3
$ws.Range("H3").Value = 918.3333
$ws.Range("I3").Value = 870.3333
$ws.Range("J3").Value = 966.3333
$ws.Range("K3").Value = 870.3333
$ws.Range("L3").Value = 966.3333
$ws.Range("M3").Value = -756.3333
$ws.Range("N3").Value = -1194.3333

# Row 134
$ws.Range("H134").Value = 27780340
$ws.Range("I134").Value = 55558180
$ws.Range("J134").Value = 2500
$ws.Range("K134").Value = 166674540
$ws.Range("L134").Value = 7500
$ws.Range("M134").Value = -166672005
$ws.Range("N134").Value = -12570

$ws = $wb.Worksheets.Item("CRP")
# Row 31
$ws.Range("H31").Value = 9561
$ws.Range("I31").Value = 1085.2941
$ws.Range("J31").Value = 30144.857
$ws.Range("K31").Value = 1085.2941
$ws.Range("L31").Value = 30144.857
$ws.Range("M31").Value = -790.2941000000001
$ws.Range("N31").Value = -30734.857

# Row 34
$ws.Range("H34").Value = 9561
$ws.Range("I34").Value = 1085.2941
$ws.Range("J34").Value = 30144.857
$ws.Range("K34").Value = 1085.2941
$ws.Range("L34").Value = 30144.857
$ws.Range("M34").Value = -883.2941000000001
$ws.Range("N34").Value = -30548.857

# Row 96
$ws.Range("H96").Value = 11712
$ws.Range("J96").Value = 11712
$ws.Range("L96").Value = 11712
$ws.Range("N96").Value = -17204

# Row 105
$ws.Range("H105").Value = 870.9167
$ws.Range("I105").Value = 826.6667
$ws.Range("J105").Value = 1003.6667
$ws.Range("K105").Value = 826.6667
$ws.Range("L105").Value = 1003.6667
$ws.Range("M105").Value = 920.3333
$ws.Range("N105").Value = -4497.6667

$ws = $wb.Worksheets.Item("CUL")
# Row 113
$ws.Range("H113").Value = 4004528.5
$ws.Range("I113").Value = 498
$ws.Range("J113").Value = 11122805
$ws.Range("K113").Value = 1494
$ws.Range("L113").Value = 33368415
$ws.Range("M113").Value = 676
$ws.Range("N113").Value = -33372755

$ws = $wb.Worksheets.Item("GSM")
# Row 122
$ws.Range("H122").Value = 83334296
$ws.Range("I122").Value = 111111610
$ws.Range("J122").Value = 2340.6667
$ws.Range("K122").Value = 333334830
$ws.Range("L122").Value = 7022.000100000001
$ws.Range("M122").Value = -333332380
$ws.Range("N122").Value = -11922.0001

# Row 126
$ws.Range("H126").Value = 1671.4
$ws.Range("I126").Value = 1181.7222
$ws.Range("J126").Value = 2930.5715
$ws.Range("K126").Value = 3545.1666
$ws.Range("L126").Value = 8791.7145
$ws.Range("M126").Value = -1075.1666
$ws.Range("N126").Value = -13731.7145

# Row 132
$ws.Range("H132").Value = 90910950
$ws.Range("I132").Value = 117647944
$ws.Range("J132").Value = 5199.2
$ws.Range("K132").Value = 352943832
$ws.Range("L132").Value = 15597.6
$ws.Range("M132").Value = -352941302
$ws.Range("N132").Value = -20657.6

$ws = $wb.Worksheets.Item("LTW")
# Row 98
$ws.Range("H98").Value = 23950
$ws.Range("J98").Value = 23950
$ws.Range("L98").Value = 23950
$ws.Range("N98").Value = -29940

# Row 136
$ws.Range("H136").Value = 2137.4146
$ws.Range("I136").Value = 2019.6562
$ws.Range("J136").Value = 2556.111
$ws.Range("K136").Value = 6058.9686
$ws.Range("L136").Value = 7668.333
$ws.Range("M136").Value = -3508.9686
$ws.Range("N136").Value = -12768.333

$ws = $wb.Worksheets.Item("WVR")
# Row 132
$ws.Range("H132").Value = 100827840
$ws.Range("I132").Value = 109092420
$ws.Range("J132").Value = 85676110
$ws.Range("K132").Value = 327277260
$ws.Range("L132").Value = 257028330
$ws.Range("M132").Value = -327274730
$ws.Range("N132").Value = -257033390

